$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 67741.01150000001
$ws.Range("C2").Value = 5416.988499999999
$ws.Range("D2").Value = 55000
$ws.Range("E2").Value = 7324.023000000008

$ws.Range("B3").Value = 64760.858
$ws.Range("C3").Value = 5210.142
$ws.Range("E3").Value = 37050.716

$ws.Range("B4").Value = 62101.3
$ws.Range("C4").Value = 5131.699999999999
$ws.Range("E4").Value = 34469.60000000001

$ws.Range("B5").Value = 61724.802
$ws.Range("C5").Value = 5049.197999999999
$ws.Range("E5").Value = 34175.60400000001

$ws.Range("B6").Value = 63342.0075
$ws.Range("C6").Value = 5084.9925
$ws.Range("E6").Value = 35757.015

$ws.Range("B7").Value = 67452.89449999999
$ws.Range("C7").Value = 5380.1055
$ws.Range("E7").Value = 39572.789

$ws.Range("B8").Value = 66444.84849999999
$ws.Range("C8").Value = 6413.1515
$ws.Range("E8").Value = 37531.69699999999

$ws.Range("B9").Value = 77477.15949999999
$ws.Range("C9").Value = 7207.8405
$ws.Range("E9").Value = 47769.31899999999

$ws.Range("B10").Value = 93062.8
$ws.Range("C10").Value = 9489.200000000001
$ws.Range("D10").Value = 55000
$ws.Range("E10").Value = 28573.60000000001

$ws.Range("B11").Value = 97938.8115
$ws.Range("C11").Value = 14318.1885
$ws.Range("E11").Value = 18620.62299999999

$ws.Range("B12").Value = 100816.808
$ws.Range("C12").Value = 15340.192
$ws.Range("E12").Value = 20476.61600000001

$ws.Range("B13").Value = 99591.12850000001
$ws.Range("C13").Value = 15589.8715
$ws.Range("E13").Value = 19001.25700000001

$ws.Range("B14").Value = 103078.6315
$ws.Range("C14").Value = 15359.3685
$ws.Range("E14").Value = 22719.26300000001

$ws.Range("B15").Value = 103456.56
$ws.Range("C15").Value = 15511.44
$ws.Range("E15").Value = 22945.12

$ws.Range("B16").Value = 106139.5625
$ws.Range("C16").Value = 15666.4375
$ws.Range("E16").Value = 25473.125

$ws.Range("B17").Value = 98840.583
$ws.Range("C17").Value = 16125.417
$ws.Range("E17").Value = 17715.166

$ws.Range("B18").Value = 0
$ws.Range("C18").Value = 16406.873
$ws.Range("E18").ClearContents()

$ws.Range("B19").Value = 91921.8355
$ws.Range("C19").Value = 16163.1645
$ws.Range("E19").Value = 10758.671

$ws.Range("B20").Value = 91677.2825
$ws.Range("C20").Value = 15659.7175
$ws.Range("E20").Value = 11017.565

$ws.Range("B21").Value = 78811.4135
$ws.Range("C21").Value = 13764.5865
$ws.Range("E21").Value = 46.8269999999975

$ws.Range("B22").Value = 80527.792
$ws.Range("C22").Value = 11506.208
$ws.Range("E22").Value = 4021.584000000003

$ws.Range("B23").Value = 59160.206
$ws.Range("C23").Value = 8102.793999999999
$ws.Range("E23").Value = -13942.588

$ws.Range("B24").Value = 65087.887
$ws.Range("C24").Value = 6308.112999999999
$ws.Range("E24").Value = -6220.225999999995

$ws.Range("C25").Value = 5711.482
$ws.Range("E25").Value = 855.6514999999927
